# Append a new data row (row 41) to Sheet1, mirroring the existing
# daily-ranking log rows (date / weekday / hour / ranking).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a date-look-alike string ("2025/09/30") that must stay
# literal text, not get auto-converted to a date serial by Excel's input
# parser. Force text via NumberFormat "@" while assigning, then clear the
# formatting again so the cell ends up with the sheet's default style,
# matching the rest of the data rows.
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = "2025/09/30"
$ws.Range("A41").ClearFormats()

$ws.Range("B41").Value = "火"
$ws.Range("C41").Value = 20
$ws.Range("D41").Value = 3
